# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# This corrects the affected team box-score stats/ranks for 2013-03-07
# and normalizes the Date column (BF) from "3-7-2012-13" to "2013-03-07".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric corrections to team stat / rank columns
$ws.Range("AD2").Value = 19
$ws.Range("AU2").Value = 2
$ws.Range("AD3").Value = 19
$ws.Range("AK4").Value = 18
$ws.Range("AP4").Value = 9
$ws.Range("AH5").Value = 16
$ws.Range("AD8").Value = 19
$ws.Range("AI8").Value = 7
$ws.Range("AP8").Value = 19
$ws.Range("AY8").Value = 4
$ws.Range("BA8").Value = 22
$ws.Range("D9").Value = 62
$ws.Range("E9").Value = 40
$ws.Range("G9").Value = 0.645
$ws.Range("J9").Value = 85.40000000000001
$ws.Range("K9").Value = 0.476
$ws.Range("N9").Value = 0.341
$ws.Range("O9").Value = 17.9
$ws.Range("P9").Value = 26
$ws.Range("Q9").Value = 0.6889999999999999
$ws.Range("R9").Value = 13.4
$ws.Range("T9").Value = 45
$ws.Range("U9").Value = 24.4
$ws.Range("Z9").Value = 20.9
$ws.Range("AA9").Value = 21.8
$ws.Range("AB9").Value = 105.7
$ws.Range("AC9").Value = 4
$ws.Range("AD9").Value = 4
$ws.Range("AF9").Value = 7
$ws.Range("AH9").Value = 9
$ws.Range("AN9").Value = 26
$ws.Range("AP9").Value = 3
$ws.Range("AR9").Value = 1
$ws.Range("AU9").Value = 3
$ws.Range("AV9").Value = 23
$ws.Range("BC9").Value = 7
$ws.Range("AD10").Value = 1
$ws.Range("AV10").Value = 22
$ws.Range("AD11").Value = 4
$ws.Range("AD12").Value = 4
$ws.Range("AN12").Value = 8
$ws.Range("AU12").Value = 6
$ws.Range("D14").Value = 63
$ws.Range("F14").Value = 19
$ws.Range("G14").Value = 0.698
$ws.Range("J14").Value = 80.7
$ws.Range("K14").Value = 0.476
$ws.Range("O14").Value = 16.6
$ws.Range("P14").Value = 23.5
$ws.Range("Q14").Value = 0.706
$ws.Range("T14").Value = 41.7
$ws.Range("W14").Value = 9.9
$ws.Range("Y14").Value = 4.2
$ws.Range("Z14").Value = 20.9
$ws.Range("AB14").Value = 100.8
$ws.Range("AC14").Value = 6.7
$ws.Range("AE14").Value = 3
$ws.Range("AF14").Value = 4
$ws.Range("AI14").Value = 6
$ws.Range("AP14").Value = 8
$ws.Range("AT14").Value = 17
$ws.Range("AU14").Value = 4
$ws.Range("AV14").Value = 16
$ws.Range("AD15").Value = 4
$ws.Range("AE16").Value = 5
$ws.Range("AI16").Value = 23
$ws.Range("AR16").Value = 2
$ws.Range("BC16").Value = 6
$ws.Range("AD19").Value = 29
$ws.Range("AP19").Value = 4
$ws.Range("AT19").Value = 10
$ws.Range("AD20").Value = 4
$ws.Range("AT20").Value = 22
$ws.Range("D21").Value = 58
$ws.Range("F21").Value = 21
$ws.Range("G21").Value = 0.638
$ws.Range("K21").Value = 0.442
$ws.Range("N21").Value = 0.37
$ws.Range("P21").Value = 21.7
$ws.Range("Q21").Value = 0.761
$ws.Range("R21").Value = 11.2
$ws.Range("S21").Value = 30.1
$ws.Range("U21").Value = 19.6
$ws.Range("Y21").Value = 3.9
$ws.Range("AB21").Value = 99.90000000000001
$ws.Range("AC21").Value = 3.8
$ws.Range("AD21").Value = 29
$ws.Range("AI21").Value = 21
$ws.Range("AJ21").Value = 14
$ws.Range("AK21").Value = 17
$ws.Range("AN21").Value = 7
$ws.Range("AP21").Value = 18
$ws.Range("AS21").Value = 21
$ws.Range("AT21").Value = 21
$ws.Range("BA21").Value = 21
$ws.Range("D22").Value = 60
$ws.Range("E22").Value = 44
$ws.Range("G22").Value = 0.733
$ws.Range("I22").Value = 38.4
$ws.Range("J22").Value = 79.59999999999999
$ws.Range("K22").Value = 0.483
$ws.Range("M22").Value = 19.5
$ws.Range("O22").Value = 22.5
$ws.Range("R22").Value = 10.4
$ws.Range("T22").Value = 43.1
$ws.Range("U22").Value = 22
$ws.Range("AB22").Value = 107
$ws.Range("AC22").Value = 9.5
$ws.Range("AD22").Value = 19
$ws.Range("AE22").Value = 3
$ws.Range("AI22").Value = 5
$ws.Range("AS22").Value = 6
$ws.Range("AT22").Value = 9
$ws.Range("AD23").Value = 4
$ws.Range("AU23").Value = 5
$ws.Range("AD24").Value = 19
$ws.Range("AT24").Value = 18
$ws.Range("AH25").Value = 16
$ws.Range("AS25").Value = 22
$ws.Range("AD26").Value = 19
$ws.Range("AN26").Value = 25
$ws.Range("AV26").Value = 15
$ws.Range("AD27").Value = 1
$ws.Range("AD28").Value = 4
$ws.Range("AD29").Value = 4
$ws.Range("AI29").Value = 22
$ws.Range("AJ29").Value = 13
$ws.Range("AK29").Value = 19
$ws.Range("AS31").Value = 5

# Date column (BF) values were stored as text in "YYYY-MM-DD" form; force text format
# so Excel does not auto-convert the string into a date serial number.
$dateCells = @("BF2","BF3","BF4","BF5","BF6","BF7","BF8","BF9","BF10","BF11","BF12","BF13","BF14","BF15","BF16","BF17","BF18","BF19","BF20","BF21","BF22","BF23","BF24","BF25","BF26","BF27","BF28","BF29","BF30","BF31")
foreach ($cell in $dateCells) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = "2013-03-07"
}
